$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 319: new positive cases corrected down by 1
$ws.Range("C319").Value = 68

# Row 323: new positive cases corrected up by 1
$ws.Range("C323").Value = 135

# Row 324: new positive cases and new hospital deaths updated
$ws.Range("C324").Value = 118
$ws.Range("L324").Value = 1

# Row 325: new positive cases, SI patients, and deaths updated
$ws.Range("C325").Value = 120
$ws.Range("G325").Value = 123
$ws.Range("L325").Value = 4
$ws.Range("M325").Value = 1

# Row 326 (2021-01-16): previously empty, now populated with data
$ws.Range("C326").Value = 56
$ws.Range("E326").Value = 11
$ws.Range("F326").Value = 9
$ws.Range("G326").Value = 116
$ws.Range("L326").Value = 1
$ws.Range("M326").Value = 2

# Row 327 (2021-01-17): previously empty, now populated with data
$ws.Range("C327").Value = 41
$ws.Range("E327").Value = 11
$ws.Range("F327").Value = 10
$ws.Range("G327").Value = 109
$ws.Range("L327").Value = 2
$ws.Range("M327").Value = 0

# Row 328 (2021-01-18): previously empty, now populated with data
$ws.Range("C328").Value = 20
$ws.Range("E328").Value = 11
$ws.Range("F328").Value = 9
$ws.Range("G328").Value = 110
$ws.Range("L328").Value = 0
$ws.Range("M328").Value = 0

$excel.CalculateFullRebuild()
